$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.335.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.32%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.868.92'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4718'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2874'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06482'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.80'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '99.61'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07794'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.871.66'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7251'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.163'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '283.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.323.33'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.9997'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007478'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.114.15'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9997'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.270'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.37%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.284'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.86'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.024'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.98'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.889'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.09657'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.318'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.95%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.484'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.228'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.133'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04802'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.123'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.52%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6888'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.77%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.713'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01895'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.841'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.48%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '75.38'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.66%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.268'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.951'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4215'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.21%  '
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8245'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.35%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '100.82'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.75%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.821'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.993'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.99'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.05761'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '884.52'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.90%  '
